$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "2C peeled chopped apple, 3/4C sugar, 2T oil, 1 egg, 1 1/2C flour. 1 1/2t  baking powder, 1/2t baking soda, 1t cinnamon, 1C cranberries, 1/2C chopped walnuts"
$ws.Range("B3").Value = "2C flour, 1C sugar, 1 1/2t baking powder, 1t salt, 1/2t baking soda, 3/4C orange juice, 1T grated orange peel, 2T shortening, 1 egg, well beaten, 1 1/2C coarsely chopped cranberries, 1/2C chopped nuts"
$ws.Range("B4").Value = "1T butter-melted, 1T boiled cider, optional, 1 1/2C whole cranberries, 1C peeled diced fresh apple, 1/4C sugar, 1 1/2C Unbleached All-Purpose Flour, heaping 1/2t salt, 1 1/2t baking powder, 1t ground cinnamon, 2/3C sugar, 2/3C applesauce (unsweetened preferred), 8Tbutter-melted, 2 large eggs, 1t vanilla extract"

$ws.Range("B5").Select() | Out-Null
